$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 1114, shifting rows 1114:1215 down to 1115:1216.
$ws.Rows.Item(1114).EntireRow.Insert()

# Populate the new row with the new weekly data point.
$ws.Range("A1114").Value = 8
$ws.Range("B1114").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1114").Value = "Coquimbo"
$ws.Range("D1114").Value = 45166
$ws.Range("D1114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E1114").Value = 4
$ws.Range("F1114").Value = 100112023
$ws.Range("G1114").Value = "Brócoli"
$ws.Range("H1114").Value = "Sin especificar"
$ws.Range("I1114").Value = "Primera"
$ws.Range("J1114").Value = 2100
$ws.Range("K1114").Value = 700
$ws.Range("L1114").Value = 800
$ws.Range("M1114").Value = 750
$ws.Range("N1114").Value = "$/unidad"
$ws.Range("O1114").Value = "Provincia del Elquí"
$ws.Range("P1114").Value = 750
$ws.Range("Q1114").Value = 1
$ws.Range("R1114").Value = "Hortaliza"
